# Rename header columns from *_old / *_new suffixes to *_FV2210 / *_FV2304,
# freeze the header row, and turn the data range into an Excel Table (ListObject).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1, columns A:U) -------------------------
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2210"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2304"
        }
    }
}

# --- 2. Freeze the header row ---------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert the data range into an Excel Table (ListObject) -----------
$dataRange = $ws.Range("A1:U72")
$listObject = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""
